$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: columns A-D, F-M become "road"; column E becomes "stoplight"
$ws.Range("A18:D18").Value = "road"
$ws.Range("E18").Value = "stoplight"
$ws.Range("F18:M18").Value = "road"

# Row 19: same pattern
$ws.Range("A19:D19").Value = "road"
$ws.Range("E19").Value = "stoplight"
$ws.Range("F19:M19").Value = "road"
